# Generate Report for Handback
#
# This reflects the localization pipeline writing back handback results:
#   - the "Ready for handoff" status becomes "Handed back: in sync with en-US"
#     everywhere it is shown (Overview + per-locale sheets)
#   - each per-locale sheet (zh-cn, de-de) gets its "Latest Target File",
#     "Latest Handback File" and "Latest Handback DateTime" columns (J/K/L)
#     populated for both rows, with J turned into a hyperlink like column A
#   - a couple of columns are widened to fit the new, longer content

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# Hyperlink font look-alike (same as the existing "HyperLink" cell style:
# single underline, cornflower blue FF6495ED == COM color 15570276)
function Set-HyperlinkLook($rng) {
    $rng.Font.Underline = 2        # xlUnderlineStyleSingle
    $rng.Font.Color = 15570276     # RGB(0x64,0x95,0xED) in BGR-long form
}

# ---------------------------------------------------------------------------
# Overview sheet: just the status text changes (columns E "zh-cn" and F
# "de-de"), for both data rows. Widen the columns to fit the longer text.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

if ($wsOverview.Range("E2").Value2 -eq $oldStatus) { $wsOverview.Range("E2").Value = $newStatus }
if ($wsOverview.Range("F2").Value2 -eq $oldStatus) { $wsOverview.Range("F2").Value = $newStatus }
if ($wsOverview.Range("E3").Value2 -eq $oldStatus) { $wsOverview.Range("E3").Value = $newStatus }
if ($wsOverview.Range("F3").Value2 -eq $oldStatus) { $wsOverview.Range("F3").Value = $newStatus }

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

if ($wsZh.Range("C2").Value2 -eq $oldStatus) { $wsZh.Range("C2").Value = $newStatus }
if ($wsZh.Range("C3").Value2 -eq $oldStatus) { $wsZh.Range("C3").Value = $newStatus }

# Row 2 - file 37d468fd...
$wsZh.Range("J2").Value = "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md"
Set-HyperlinkLook($wsZh.Range("J2"))
$wsZh.Hyperlinks.Add($wsZh.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d931f1880d4e14e688631da00d1895efcb611d17/e2e/37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md", "", "", "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md")
$wsZh.Range("K2").Value = "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.d56cd217d28f87fab3f4caaa3a30e3842726b105.zh-cn.xlf"
$wsZh.Range("L2").Value = "2016-12-16 10:19:18"

# Row 3 - file 8593562f...
$wsZh.Range("J3").Value = "8593562f-b597-4f5e-8669-a974143ad71e.md"
Set-HyperlinkLook($wsZh.Range("J3"))
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d931f1880d4e14e688631da00d1895efcb611d17/e2e/8593562f-b597-4f5e-8669-a974143ad71e.md", "", "", "8593562f-b597-4f5e-8669-a974143ad71e.md")
$wsZh.Range("K3").Value = "8593562f-b597-4f5e-8669-a974143ad71e.7f70fcb6fa205f2832122286ee0bdc78f0c6c6ad.zh-cn.xlf"
$wsZh.Range("L3").Value = "2016-12-16 10:19:18"

$wsZh.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZh.Columns.Item(10).ColumnWidth = 40
$wsZh.Columns.Item(11).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

if ($wsDe.Range("C2").Value2 -eq $oldStatus) { $wsDe.Range("C2").Value = $newStatus }
if ($wsDe.Range("C3").Value2 -eq $oldStatus) { $wsDe.Range("C3").Value = $newStatus }

# Row 2 - file 37d468fd...
$wsDe.Range("J2").Value = "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md"
Set-HyperlinkLook($wsDe.Range("J2"))
$wsDe.Hyperlinks.Add($wsDe.Range("J2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d931f1880d4e14e688631da00d1895efcb611d17/e2e/37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md", "", "", "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.md")
$wsDe.Range("K2").Value = "37d468fd-95b3-48c7-abfe-c4d62c53ff3e.d56cd217d28f87fab3f4caaa3a30e3842726b105.de-de.xlf"
$wsDe.Range("L2").Value = "2016-12-16 10:19:36"

# Row 3 - file 8593562f...
$wsDe.Range("J3").Value = "8593562f-b597-4f5e-8669-a974143ad71e.md"
Set-HyperlinkLook($wsDe.Range("J3"))
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d931f1880d4e14e688631da00d1895efcb611d17/e2e/8593562f-b597-4f5e-8669-a974143ad71e.md", "", "", "8593562f-b597-4f5e-8669-a974143ad71e.md")
$wsDe.Range("K3").Value = "8593562f-b597-4f5e-8669-a974143ad71e.7f70fcb6fa205f2832122286ee0bdc78f0c6c6ad.de-de.xlf"
$wsDe.Range("L3").Value = "2016-12-16 10:19:36"

$wsDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDe.Columns.Item(10).ColumnWidth = 40
$wsDe.Columns.Item(11).ColumnWidth = 40
